$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102, shifting existing rows 102:177 down to 103:178.
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with its data (same static columns as the
# surrounding "Femacal de La Calera" / "Achicoria" records, new measurement values).
$ws.Cells.Item(102, 1).Value = 3
$ws.Cells.Item(102, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(102, 3).Value = "Coquimbo"
$ws.Cells.Item(102, 4).Value = 44574
$ws.Cells.Item(102, 5).Value = 5
$ws.Cells.Item(102, 6).Value = 100112010
$ws.Cells.Item(102, 7).Value = "Achicoria"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 108
$ws.Cells.Item(102, 11).Value = 5500
$ws.Cells.Item(102, 12).Value = 6000
$ws.Cells.Item(102, 13).Value = 5769
$ws.Cells.Item(102, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(102, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(102, 16).Value = 361
$ws.Cells.Item(102, 17).Value = 16
$ws.Cells.Item(102, 18).Value = "Hortaliza"
